# Delete the data row for account 005171652 / Bruno / 200 (Excel row 71)
# from the "Export" sheet, shifting all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$ws.Rows.Item(71).Delete()
